# Auto-generated edit script applying value changes per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AA2").Value = 38
$ws.Range("AB2").Value = 18
$ws.Range("AC2").Value = 9.4
$ws.Range("AD2").Value = 12.5
$ws.Range("AF2").Value = 25
$ws.Range("AH2").Value = 14.5
$ws.Range("AI2").Value = 29
$ws.Range("AJ2").Value = 44
$ws.Range("AK2").Value = 25
$ws.Range("AL2").Value = 30
$ws.Range("AM2").Value = 55
$ws.Range("AN2").Value = 16
$ws.Range("AO2").Value = 13
$ws.Range("F2").Value = 2.8
$ws.Range("H2").Value = 2.5
$ws.Range("J2").Value = 3.85
$ws.Range("N2").Value = 5.9
$ws.Range("O2").Value = 1.18
$ws.Range("P2").Value = 2.64
$ws.Range("S2").Value = 2.42
$ws.Range("T2").Value = 1.52
$ws.Range("U2").Value = 2.74
$ws.Range("V2").Value = 1.62
$ws.Range("X2").Value = 25
$ws.Range("Y2").Value = 17
$ws.Range("Z2").Value = 21
$ws.Range("AB3").Value = 1000
$ws.Range("AC3").Value = 1000
$ws.Range("AD3").Value = 1000
$ws.Range("AH3").Value = 1000
$ws.Range("AL3").Value = 21
$ws.Range("AN3").Value = 4.6
$ws.Range("F3").Value = 1.5
$ws.Range("G3").Value = 1.58
$ws.Range("H3").Value = 5
$ws.Range("K3").Value = 6
$ws.Range("L3").Value = 1.19
$ws.Range("P3").Value = 3.2
$ws.Range("W3").Value = 2.72
$ws.Range("AB4").Value = 7.8
$ws.Range("AF4").Value = 9.6
$ws.Range("AJ4").Value = 16
$ws.Range("AL4").Value = 48
$ws.Range("F4").Value = 1.53
$ws.Range("G4").Value = 1.57
$ws.Range("H4").Value = 7.6
$ws.Range("I4").Value = 8.199999999999999
$ws.Range("K4").Value = 4.6
$ws.Range("O4").Value = 1.32
$ws.Range("P4").Value = 1.9
$ws.Range("Q4").Value = 1.94
$ws.Range("W4").Value = 2.76
$ws.Range("F5").Value = 2.78
$ws.Range("H5").Value = 2.62
$ws.Range("I5").Value = 3.1
$ws.Range("J5").Value = 3
$ws.Range("K5").Value = 3.4
$ws.Range("N5").Value = 2.74
$ws.Range("O5").Value = 1.43
$ws.Range("P5").Value = 1.58
$ws.Range("Q5").Value = 2.3
$ws.Range("R5").Value = 1.21
$ws.Range("S5").Value = 4.1
$ws.Range("T5").Value = 1.92
$ws.Range("U5").Value = 1.85
$ws.Range("V5").Value = 1.47
$ws.Range("W5").Value = 1.44
$ws.Range("G6").Value = 3.85
$ws.Range("H6").Value = 2.16
$ws.Range("I6").Value = 2.3
$ws.Range("J6").Value = 3.45
$ws.Range("P6").Value = 1.92
$ws.Range("Q6").Value = 1.95
$ws.Range("S6").Value = 3.55
$ws.Range("V6").Value = 1.76
$ws.Range("AK7").Value = 18
$ws.Range("F7").Value = 1.32
$ws.Range("S7").Value = 2.62
$ws.Range("T7").Value = 2.14
$ws.Range("V7").Value = 1.08
$ws.Range("H8").Value = 2.58
$ws.Range("I8").Value = 2.66
$ws.Range("N8").Value = 3.8
$ws.Range("P8").Value = 1.98
$ws.Range("W8").Value = 1.5
$ws.Range("AA9").Value = 24
$ws.Range("AB9").Value = 23
$ws.Range("AF9").Value = 65
$ws.Range("AG9").Value = 18.5
$ws.Range("AN9").Value = 40
$ws.Range("G9").Value = 4.2
$ws.Range("H9").Value = 1.9
$ws.Range("J9").Value = 4.1
$ws.Range("K9").Value = 4.2
$ws.Range("N9").Value = 5.5
$ws.Range("X9").Value = 40
$ws.Range("AO10").Value = 50
$ws.Range("F10").Value = 2.3
$ws.Range("H10").Value = 3.3
$ws.Range("I10").Value = 3.7
$ws.Range("K10").Value = 3.6
$ws.Range("O10").Value = 1.33
$ws.Range("P10").Value = 1.85
$ws.Range("Q10").Value = 1.99
$ws.Range("S10").Value = 3.55
$ws.Range("T10").Value = 1.76
$ws.Range("M11").Value = 1.01
$ws.Range("AA12").Value = 44
$ws.Range("AE12").Value = 28
$ws.Range("AI12").Value = 1000
$ws.Range("AO12").Value = 15.5
$ws.Range("F12").Value = 2.7
$ws.Range("G12").Value = 2.98
$ws.Range("H12").Value = 2.4
$ws.Range("I12").Value = 2.58
$ws.Range("J12").Value = 3.95
$ws.Range("L12").Value = 1.29
$ws.Range("N12").Value = 5.5
$ws.Range("O12").Value = 1.19
$ws.Range("P12").Value = 2.56
$ws.Range("U12").Value = 2.62
$ws.Range("V12").Value = 1.63
$ws.Range("W12").Value = 1.5
$ws.Range("Y12").Value = 19.5
$ws.Range("F13").Value = 2.82
$ws.Range("I13").Value = 2.84
$ws.Range("N13").Value = 3.6
$ws.Range("AO15").Value = 29
$ws.Range("G15").Value = 4.5
$ws.Range("I15").Value = 2.36
$ws.Range("J15").Value = 3.15
$ws.Range("K15").Value = 3.55
$ws.Range("R15").Value = 1.24
$ws.Range("S15").Value = 4.4
$ws.Range("T15").Value = 1.83
$ws.Range("V15").Value = 1.73
$ws.Range("W15").Value = 1.31
$ws.Range("AA16").Value = 95
$ws.Range("AB16").Value = 8.199999999999999
$ws.Range("AC16").Value = 7.2
$ws.Range("AE16").Value = 200
$ws.Range("AG16").Value = 19.5
$ws.Range("AO16").Value = 70
$ws.Range("F16").Value = 3.45
$ws.Range("G16").Value = 3.75
$ws.Range("H16").Value = 2.56
$ws.Range("I16").Value = 2.66
$ws.Range("L16").Value = 1.75
$ws.Range("M16").Value = 1.19
$ws.Range("T16").Value = 2.58
$ws.Range("V16").Value = 1.6
$ws.Range("W16").Value = 1.36
$ws.Range("Y16").Value = 6.6
$ws.Range("Z16").Value = 14.5
$ws.Range("AH17").Value = 1000
$ws.Range("F17").Value = 2.02
$ws.Range("G17").Value = 2.12
$ws.Range("H17").Value = 4.5
$ws.Range("J17").Value = 3.05
$ws.Range("K17").Value = 3.4
$ws.Range("N17").Value = 2.5
$ws.Range("P17").Value = 1.5
$ws.Range("S17").Value = 5.5
$ws.Range("T17").Value = 2.22
$ws.Range("U17").Value = 1.69
$ws.Range("V17").Value = 1.25
$ws.Range("W17").Value = 1.9
$ws.Range("F18").Value = 4.6
$ws.Range("G18").Value = 5.7
$ws.Range("H18").Value = 1.87
$ws.Range("I18").Value = 2.06
$ws.Range("K18").Value = 3.65
$ws.Range("N18").Value = 2.48
$ws.Range("O18").Value = 1.54
$ws.Range("V18").Value = 1.94
$ws.Range("AA19").Value = 140
$ws.Range("AC19").Value = 8.199999999999999
$ws.Range("AD19").Value = 980
$ws.Range("AE19").Value = 80
$ws.Range("AF19").Value = 11.5
$ws.Range("AI19").Value = 95
$ws.Range("AJ19").Value = 26
$ws.Range("AK19").Value = 27
$ws.Range("AL19").Value = 55
$ws.Range("AM19").Value = 160
$ws.Range("AN19").Value = 18.5
$ws.Range("F19").Value = 1.87
$ws.Range("G19").Value = 1.97
$ws.Range("H19").Value = 4.8
$ws.Range("I19").Value = 5.4
$ws.Range("J19").Value = 3.55
$ws.Range("K19").Value = 3.7
$ws.Range("L19").Value = 1.46
$ws.Range("M19").Value = 1.09
$ws.Range("N19").Value = 3.2
$ws.Range("O19").Value = 1.39
$ws.Range("P19").Value = 1.76
$ws.Range("Q19").Value = 2.14
$ws.Range("R19").Value = 1.28
$ws.Range("T19").Value = 1.94
$ws.Range("V19").Value = 1.22
$ws.Range("W19").Value = 2.04
$ws.Range("X19").Value = 13.5
$ws.Range("Y19").Value = 16.5
$ws.Range("AA20").Value = 980
$ws.Range("AB20").Value = 980
$ws.Range("AC20").Value = 980
$ws.Range("AD20").Value = 980
$ws.Range("AE20").Value = 980
$ws.Range("AF20").Value = 980
$ws.Range("AG20").Value = 980
$ws.Range("AH20").Value = 980
$ws.Range("AO20").Value = 980
$ws.Range("F20").Value = 3.65
$ws.Range("I20").Value = 2.38
$ws.Range("N20").Value = 2.94
$ws.Range("S20").Value = 4.4
$ws.Range("X20").Value = 980
$ws.Range("Y20").Value = 980
$ws.Range("Z20").Value = 980
